$p = $ppt.ActivePresentation

# Locate the slide that holds the "Golge Degisken Tuzagi" (dummy-variable
# trap) example -- this is the slide whose PowerPoint-internal SlideID is
# 271 (creationId 3773428193 in the change log).
$s = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $cand = $p.Slides.Item($i)
    if ($cand.SlideID -eq 271) {
        $s = $cand
    }
}
if ($s -eq $null) {
    $s = $p.Slides.Item(8)
}

# --- Update the two lookup-table headers ---
# "Age" -> "Ulke", "Salary" -> "Yas", "Purchased" -> "Gelir"
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTable) {
        $tbl = $sh.Table
        $tbl.Cell(1,1).Shape.TextFrame.TextRange.Text = "Ulke"
        $tbl.Cell(1,2).Shape.TextFrame.TextRange.Text = "Yas"
        $tbl.Cell(1,3).Shape.TextFrame.TextRange.Text = "Gelir"
    }
}

# --- Remove the click-to-reveal animation that targeted the rectangle below ---
$seq = $s.TimeLine.MainSequence
for ($i = $seq.Count; $i -ge 1; $i--) {
    $eff = $seq.Item($i)
    if ($eff.Shape.Name -eq "Dikdörtgen 17") {
        $eff.Delete()
    }
}

# --- Remove the now-unneeded rectangle shape (id 18) itself ---
for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Dikdörtgen 17") {
        $sh.Delete()
    }
}
